$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (price) values stay as text, matching the source data format
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '24.707.63'
$ws.Range('E2').Value = '  +2.38%  '

$ws.Range('D3').Value = '1.703.55'
$ws.Range('E3').Value = '  +1.77%  '

$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  -0.39%  '

$ws.Range('D5').Value = '308.87'
$ws.Range('E5').Value = '  +0.20%  '

$ws.Range('D6').Value = '0.9979'
$ws.Range('E6').Value = '  -0.31%  '

$ws.Range('D7').Value = '0.3741'
$ws.Range('E7').Value = '  +0.52%  '

$ws.Range('D8').Value = '49.25'
$ws.Range('E8').Value = '  +3.70%  '

$ws.Range('D9').Value = '0.3439'
$ws.Range('E9').Value = '  +0.17%  '

$ws.Range('D10').Value = '1.188'
$ws.Range('E10').Value = '  +0.17%  '

$ws.Range('E11').Value = '  +2.39%  '

$ws.Range('D12').Value = '0.9992'
$ws.Range('E12').Value = '  -0.51%  '

$ws.Range('D13').Value = '20.90'
$ws.Range('E13').Value = '  +2.67%  '

$ws.Range('D14').Value = '6.237'
$ws.Range('E14').Value = '  +2.55%  '

$ws.Range('D15').Value = '6.954'
$ws.Range('E15').Value = '  +3.01%  '

$ws.Range('D16').Value = '1.704.84'
$ws.Range('E16').Value = '  +1.43%  '

$ws.Range('E17').Value = '  +1.73%  '

$ws.Range('D18').Value = '0.06714'
$ws.Range('E18').Value = '  +0.18%  '

$ws.Range('D19').Value = '0.9982'
$ws.Range('E19').Value = '  -0.30%  '

$ws.Range('D20').Value = '84.06'
$ws.Range('E20').Value = '  +3.07%  '

$ws.Range('E21').Value = '  +4.17%  '

$ws.Range('D22').Value = '6.322'
$ws.Range('E22').Value = '  +3.81%  '

$ws.Range('D23').Value = '13.04'
$ws.Range('E23').Value = '  +8.83%  '

$ws.Range('D24').Value = '24.708.19'
$ws.Range('E24').Value = '  +2.41%  '

$ws.Range('D25').Value = '2.424'
$ws.Range('E25').Value = '  +0.35%  '

$ws.Range('D26').Value = '2.766'
$ws.Range('E26').Value = '  +3.95%  '

$ws.Range('E27').Value = '  +3.12%  '

$ws.Range('D28').Value = '150.29'
$ws.Range('E28').Value = '  -1.93%  '

$ws.Range('D29').Value = '131.06'
$ws.Range('E29').Value = '  +3.41%  '

$ws.Range('D30').Value = '1.893.34'
$ws.Range('E30').Value = '  +1.41%  '

$ws.Range('D31').Value = '1.182'
$ws.Range('E31').Value = '  +20.85%  '

$ws.Range('D32').Value = '6.764'
$ws.Range('E32').Value = '  +6.32%  '

$ws.Range('D33').Value = '4.169'
$ws.Range('E33').Value = '  +1.47%  '

$ws.Range('D34').Value = '1.798'
$ws.Range('E34').Value = '  +1.36%  '

$ws.Range('D35').Value = '0.08853'
$ws.Range('E35').Value = '  +4.74%  '

$ws.Range('D36').Value = '13.65'
$ws.Range('E36').Value = '  +10.73%  '

$ws.Range('D37').Value = '5.518'
$ws.Range('E37').Value = '  +3.38%  '

$ws.Range('D38').Value = '0.06566'
$ws.Range('E38').Value = '  +1.96%  '

$ws.Range('D39').Value = '8.954'
$ws.Range('E39').Value = '  +0.78%  '

$ws.Range('D40').Value = '0.02378'
$ws.Range('E40').Value = '  +2.01%  '

$ws.Range('D41').Value = '0.2226'
$ws.Range('E41').Value = '  +5.07%  '

$ws.Range('D42').Value = '1.275'
$ws.Range('E42').Value = '  +0.69%  '

$ws.Range('D43').Value = '0.6418'
$ws.Range('E43').Value = '  +4.25%  '

$ws.Range('D44').Value = '0.9977'
$ws.Range('E44').Value = '  -0.24%  '

$ws.Range('D45').Value = '13.93'
$ws.Range('E45').Value = '  +6.44%  '

$ws.Range('D46').Value = '0.6104'
$ws.Range('E46').Value = '  +3.02%  '

$ws.Range('D47').Value = '3.807'
$ws.Range('E47').Value = '  +0.21%  '

$ws.Range('D48').Value = '2.110'
$ws.Range('E48').Value = '  +4.09%  '

$ws.Range('D49').Value = '129.63'
$ws.Range('E49').Value = '  +2.57%  '

$ws.Range('D50').Value = '0.07287'
$ws.Range('E50').Value = '  +1.58%  '

$ws.Range('D51').Value = '79.14'
$ws.Range('E51').Value = '  +4.22%  '
